$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.460.40"
$ws.Range("E2").Value = "  +2.42%  "

$ws.Range("D3").Value = "4.061.85"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("D7").Value = "4.052.02"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.697"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.87%  "

$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.12%  "

$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.78%  "

$ws.Range("D15").Value = "4.709.09"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").Value = "4.060.21"
$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.75%  "

$ws.Range("E19").Value = "  +4.30%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").Value = "73.368.93"
$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.52%  "

$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.69%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "98.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +20.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.82%  "

$ws.Range("E33").Value = "  +5.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "686.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "68.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.38%  "

$ws.Range("D38").Value = "0.0₃0916"
$ws.Range("E38").Value = "  +10.57%  "

$ws.Range("E39").Value = "  +6.66%  "

$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("E41").Value = "  -1.45%  "

$ws.Range("E42").Value = "  +18.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("E45").Value = "  +3.61%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("E47").Value = "  +2.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.67%  "

$ws.Range("E49").Value = "  +8.92%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.48%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
